# The deck had two "Preliminary Planning" slides back-to-back at positions
# 9 and 10 (position 9 holding the "Program Flowchart" diagram). This move
# relocates the flowchart slide (originally at position 9) so it lands
# after the "4. Reflection / Refactor Continued" slide (originally at
# position 18), shifting the slides in between up by one position.
$p = $ppt.ActivePresentation
$p.Slides.Range(9).MoveTo(18)
